$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.375.16'
$ws.Range('E2').Value = '  +4.30%  '

$ws.Range('D3').Value = '3.616.66'
$ws.Range('E3').Value = '  +4.53%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '589.72'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.17%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '191.19'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.94%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.644'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.44%  '

$ws.Range('D8').Value = '3.610.39'
$ws.Range('E8').Value = '  +4.52%  '

$ws.Range('E9').Value = '  -0.09%  '

$ws.Range('E10').Value = '  -0.22%  '

$ws.Range('E11').Value = '  +2.26%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '58.23'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +4.91%  '

$ws.Range('E13').Value = '  +3.10%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.80'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.13%  '

$ws.Range('D15').Value = '4.191.57'
$ws.Range('E15').Value = '  +4.83%  '

$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.42'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +4.44%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.611.01'
$ws.Range('E17').Value = '  +4.52%  '

$ws.Range('D18').Value = '70.297.07'
$ws.Range('E18').Value = '  +4.33%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.50'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.63%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.121'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.29%  '

$ws.Range('E21').Value = '  +3.87%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '493.61'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.66%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.53'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +16.53%  '

$ws.Range('E24').Value = '  +8.53%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.46'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +6.43%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '90.89'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.06%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.10'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.88%  '

$ws.Range('E28').Value = '  +0.75%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.46'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +5.52%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.48'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.83%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.52'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +7.82%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '629.52'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +6.02%  '

$ws.Range('E33').Value = '  +5.14%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.118'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +6.94%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '65.10'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.97%  '

$ws.Range('D36').Value = '0.0₃0822'
$ws.Range('E36').Value = '  +4.13%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '38.15'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +4.12%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.405'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +3.54%  '

$ws.Range('E39').Value = '  +0.09%  '

$ws.Range('E40').Value = '  -1.85%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.63'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.53%  '

$ws.Range('D42').Value = '3.300.04'
$ws.Range('E42').Value = '  +4.96%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.10'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +5.43%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0446'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +4.58%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.67'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.55%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.31'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.47%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.138'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.69%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.13'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +4.11%  '

$ws.Range('E49').Value = '  -3.69%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.33'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +6.66%  '

$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '143.01'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.46%  '
